$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.935.19"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.038.78"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.657"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.64"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0769"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.877"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.82%  "
$ws.Range("D14").Value = "2.335.52"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "2.035.89"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.08%  "
$ws.Range("D18").Value = "36.910.39"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "235.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -3.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.14%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0612"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0864"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("E36").Value = "  +6.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  -12.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.39%  "
$ws.Range("D46").Value = "1.293.62"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("E49").Value = "  +8.72%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "2.222.10"
$ws.Range("E51").Value = "  -0.88%  "
